# Hortaliza / Agrícola del Norte S.A. de Arica - Pimiento
# Weekly update: insert 3 new "Zafiro rojo" price rows (fecha 44706) right
# after the existing 44664 block (original row 573), pushing the rest of
# the table down by 3 rows. The oldest 3 rows that fall off the bottom of
# the historical window are appended again at the end (rows 623-625),
# exactly mirroring the previous 3 rows (fecha 44376, Zafiro verde).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 blank rows before row 573 - shifts old rows 573:622 down to 576:625
$ws.Rows("573:575").Insert()

# 2) New rows 573-575 share the same "constant" columns as the rest of the
#    table for this market/category (copy them from row 576, which now
#    holds the data that used to be in row 573 before the insert).
$templateRow = 576
for ($i = 0; $i -lt 3; $i++) {
    $r = 573 + $i
    $ws.Cells.Item($r, 1).Value  = $ws.Cells.Item($templateRow, 1).Value()   # A Mercado ID
    $ws.Cells.Item($r, 2).Value  = $ws.Cells.Item($templateRow, 2).Value()   # B Mercado
    $ws.Cells.Item($r, 3).Value  = $ws.Cells.Item($templateRow, 3).Value()   # C Región
    $ws.Cells.Item($r, 5).Value  = $ws.Cells.Item($templateRow, 5).Value()   # E Codreg
    $ws.Cells.Item($r, 6).Value  = $ws.Cells.Item($templateRow, 6).Value()   # F Categoría ID
    $ws.Cells.Item($r, 7).Value  = $ws.Cells.Item($templateRow, 7).Value()   # G Categoría
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($templateRow, 14).Value()  # N Unidad de comercialización
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($templateRow, 15).Value()  # O Origen
    $ws.Cells.Item($r, 17).Value = $ws.Cells.Item($templateRow, 17).Value()  # Q Kg o Unidades
    $ws.Cells.Item($r, 18).Value = $ws.Cells.Item($templateRow, 18).Value()  # R Clasificación
}

# 3) Row-specific values for the 3 new rows (D, H, I, J, K, L, M, P)
$ws.Range("D573").Value = 44706
$ws.Range("H573").Value = "Zafiro rojo"
$ws.Range("I573").Value = "Primera"
$ws.Range("J573").Value = 120
$ws.Range("K573").Value = 35000
$ws.Range("L573").Value = 37000
$ws.Range("M573").Value = 36000
$ws.Range("P573").Value = 2400

$ws.Range("D574").Value = 44706
$ws.Range("H574").Value = "Zafiro rojo"
$ws.Range("I574").Value = "Segunda"
$ws.Range("J574").Value = 160
$ws.Range("K574").Value = 31000
$ws.Range("L574").Value = 33000
$ws.Range("M574").Value = 32000
$ws.Range("P574").Value = 2133

$ws.Range("D575").Value = 44706
$ws.Range("H575").Value = "Zafiro rojo"
$ws.Range("I575").Value = "Tercera"
$ws.Range("J575").Value = 160
$ws.Range("K575").Value = 27000
$ws.Range("L575").Value = 28000
$ws.Range("M575").Value = 27500
$ws.Range("P575").Value = 1833
